$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "308.49"
$r.Style = "Normal"
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "-4.38%"
$r.Style = "Normal"
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "48.52"
$r.Style = "Normal"
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "-2.35%"
$r.Style = "Normal"
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "5.184"
$r.Style = "Normal"
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "-3.38%"
$r.Style = "Normal"
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "0.07750"
$r.Style = "Normal"
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "-4.82%"
$r.Style = "Normal"
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "4.491"
$r.Style = "Normal"
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "-2.42%"
$r.Style = "Normal"
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "1.339"
$r.Style = "Normal"
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "15.12%"
$r.Style = "Normal"
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "1.558"
$r.Style = "Normal"
$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "-6.87%"
$r.Style = "Normal"
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.1225"
$r.Style = "Normal"
$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = "-9.45%"
$r.Style = "Normal"
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.1941"
$r.Style = "Normal"
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "-1.66%"
$r.Style = "Normal"
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.09390"
$r.Style = "Normal"
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "-2.08%"
$r.Style = "Normal"
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.04635"
$r.Style = "Normal"
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "1.12%"
$r.Style = "Normal"
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "0.1046"
$r.Style = "Normal"
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "-0.13%"
$r.Style = "Normal"
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "0.001260"
$r.Style = "Normal"
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "-5.38%"
$r.Style = "Normal"
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "-3.02%"
$r.Style = "Normal"
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "0.005859"
$r.Style = "Normal"
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "0.44%"
$r.Style = "Normal"
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "3.328"
$r.Style = "Normal"
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "-1.69%"
$r.Style = "Normal"
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "-6.75%"
$r.Style = "Normal"
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "0.3487"
$r.Style = "Normal"
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "2.75%"
$r.Style = "Normal"
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "8.352"
$r.Style = "Normal"
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "2.89%"
$r.Style = "Normal"
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "0.1338"
$r.Style = "Normal"
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "-5.06%"
$r.Style = "Normal"
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "-0.44%"
$r.Style = "Normal"
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "-2.35%"
$r.Style = "Normal"
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "0.004176"
$r.Style = "Normal"
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "-3.06%"
$r.Style = "Normal"
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "0.26%"
$r.Style = "Normal"
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "-4.04%"
$r.Style = "Normal"
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.02564"
$r.Style = "Normal"
$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = "-7.18%"
$r.Style = "Normal"
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.05862"
$r.Style = "Normal"
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = "6.02%"
$r.Style = "Normal"
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "73.65%"
$r.Style = "Normal"
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.007914"
$r.Style = "Normal"
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "1.08%"
$r.Style = "Normal"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.1419"
$r.Style = "Normal"
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "-2.01%"
$r.Style = "Normal"
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.008348"
$r.Style = "Normal"
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "8.78%"
$r.Style = "Normal"
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.007707"
$r.Style = "Normal"
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "-4.67%"
$r.Style = "Normal"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.3099"
$r.Style = "Normal"
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "-11.43%"
$r.Style = "Normal"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.00006963"
$r.Style = "Normal"
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = "2.96%"
$r.Style = "Normal"
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.00000000751"
$r.Style = "Normal"
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "0.09%"
$r.Style = "Normal"
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = "0.15%"
$r.Style = "Normal"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.00002102"
$r.Style = "Normal"
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = "0.09%"
$r.Style = "Normal"
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.0002002"
$r.Style = "Normal"
$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = "0.09%"
$r.Style = "Normal"
